$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Part 1"

# Cells that will hold numeric-looking text (must stay stored as text, like the
# rest of this workbook's cells, rather than being auto-coerced to numbers).
# NumberFormat must be set to Text ("@") *before* the value is assigned.
$textCellRefs = @("B2", "B3", "B4", "B7", "B8")
foreach ($ref in $textCellRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Row 1 unchanged
$ws.Range("A1").Value = "Inputs"
$ws.Range("B1").Value = "Values"

# Row 2: new "Desired P(block)" row
$ws.Range("A2").Value = "Desired P(block)"
$ws.Range("B2").Value = "0.2"

# Row 3: Arrival Rate (value changed to 5.0)
$ws.Range("A3").Value = "Arrival Rate"
$ws.Range("B3").Value = "5.0"

# Row 4: Service Rate (value changed to 6.0)
$ws.Range("A4").Value = "Service Rate"
$ws.Range("B4").Value = "6.0"

# Row 5: now empty (old "E(W) Less Than" row removed entirely)
$ws.Range("A5:B5").ClearContents()

# Row 6: Results stays the same; B6 must stay empty
$ws.Range("A6").Value = "Results"
$ws.Range("B6").ClearContents()

# Row 7: Number of Servers (value changed to 2)
$ws.Range("A7").Value = "Number of Servers"
$ws.Range("B7").Value = "2"

# Row 8: new "Actual P(block)" row
$ws.Range("A8").Value = "Actual P(block)"
$ws.Range("B8").Value = "0.1592356687898089"

# Restore default styling on the cells we had to coerce to text, so the sheet
# doesn't carry a residual explicit-style index.
foreach ($ref in $textCellRefs) {
    $ws.Range($ref).Style = "Normal"
}
